{"js": "/*\n * Replace the text of each arithmetic-equation table cell (in document\n * order: row-major, left-to-right, top-to-bottom) with its new value.\n * The table has a fixed 5-column layout; old values are read from the\n * document to preserve formatting while the text is swapped in bulk via\n * Table.values, which keeps each cell's existing run formatting intact.\n */\nconst newEquations = [\"2+14=16\", \"68-44=24\", \"61+28=89\", \"22+20=42\", \"60+26=86\", \"65+6=71\", \"15+23=38\", \"9-3=6\", \"31+12=43\", \"18+25=43\", \"26+34=60\", \"94-38=56\", \"34+21=55\", \"8+56=64\", \"76-31=45\", \"98-30=68\", \"65-48=17\", \"18+67=85\", \"92-51=41\", \"70-3=67\", \"44+20=64\", \"13-2=11\", \"66-28=38\", \"37-0=37\", \"53-4=49\", \"39+25=64\", \"14+22=36\", \"70-43=27\", \"97-24=73\", \"18+55=73\", \"79-44=35\", \"66-25=41\", \"6+81=87\", \"70-40=30\", \"10+70=80\", \"36+56=92\", \"65-44=21\", \"64-56=8\", \"9+79=88\", \"56+29=85\", \"27+18=45\", \"58+14=72\", \"58+40=98\", \"39-34=5\", \"45+52=97\", \"4+7=11\", \"55-23=32\", \"14+31=45\", \"54+22=76\", \"62+16=78\", \"68+12=80\", \"2+76=78\", \"82-2=80\", \"66+18=84\", \"70+29=99\", \"26+12=38\", \"87-58=29\", \"68-58=10\", \"31+24=55\", \"5+19=24\", \"51-28=23\", \"31+49=80\", \"58-4=54\", \"36+20=56\", \"21+59=80\", \"38+7=45\", \"61-18=43\", \"79-4=75\", \"37+47=84\", \"7+37=44\", \"67-63=4\", \"93+4=97\", \"7+4=11\", \"53-51=2\", \"0+38=38\", \"9-1=8\", \"99-93=6\", \"63-34=29\", \"74-16=58\", \"18+27=45\", \"39+59=98\", \"37+25=62\", \"40+23=63\", \"53-2=51\", \"62-29=33\", \"73-8=65\", \"60-25=35\", \"0+1=1\", \"64+17=81\", \"10+59=69\", \"57-31=26\", \"38-0=38\", \"42+11=53\", \"52-28=24\", \"2+63=65\", \"77+21=98\", \"27-17=10\", \"76-45=31\", \"76-38=38\", \"77-41=36\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst colCount = oldValues.length > 0 ? oldValues[0].length : 0;\n\nif (table.rowCount * colCount !== newEquations.length) {\n  throw new Error(\n    `Table shape (${table.rowCount}x${colCount}) does not match the ` +\n    `${newEquations.length} replacement values.`\n  );\n}\n\n// Re-shape the flat, in-order list of new values into a 2D grid matching\n// the table's existing row/column layout, then assign it back in one\n// shot so each cell's existing run/paragraph formatting is preserved\n// (only the <w:t> text content changes).\nconst newValues = [];\nfor (let r = 0; r < table.rowCount; r++) {\n  const row = [];\n  for (let c = 0; c < colCount; c++) {\n    row.push(newEquations[r * colCount + c]);\n  }\n  newValues.push(row);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the text of each arithmetic-equation table cell (in document\n# order: row-major, left-to-right, top-to-bottom) with its new value,\n# preserving each cell's existing run/paragraph formatting by assigning\n# directly to Cell(r, c).Range.Text (does not touch rPr/pPr).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"2+14=16\", \"68-44=24\", \"61+28=89\", \"22+20=42\", \"60+26=86\"),\n    @(\"65+6=71\", \"15+23=38\", \"9-3=6\", \"31+12=43\", \"18+25=43\"),\n    @(\"26+34=60\", \"94-38=56\", \"34+21=55\", \"8+56=64\", \"76-31=45\"),\n    @(\"98-30=68\", \"65-48=17\", \"18+67=85\", \"92-51=41\", \"70-3=67\"),\n    @(\"44+20=64\", \"13-2=11\", \"66-28=38\", \"37-0=37\", \"53-4=49\"),\n    @(\"39+25=64\", \"14+22=36\", \"70-43=27\", \"97-24=73\", \"18+55=73\"),\n    @(\"79-44=35\", \"66-25=41\", \"6+81=87\", \"70-40=30\", \"10+70=80\"),\n    @(\"36+56=92\", \"65-44=21\", \"64-56=8\", \"9+79=88\", \"56+29=85\"),\n    @(\"27+18=45\", \"58+14=72\", \"58+40=98\", \"39-34=5\", \"45+52=97\"),\n    @(\"4+7=11\", \"55-23=32\", \"14+31=45\", \"54+22=76\", \"62+16=78\"),\n    @(\"68+12=80\", \"2+76=78\", \"82-2=80\", \"66+18=84\", \"70+29=99\"),\n    @(\"26+12=38\", \"87-58=29\", \"68-58=10\", \"31+24=55\", \"5+19=24\"),\n    @(\"51-28=23\", \"31+49=80\", \"58-4=54\", \"36+20=56\", \"21+59=80\"),\n    @(\"38+7=45\", \"61-18=43\", \"79-4=75\", \"37+47=84\", \"7+37=44\"),\n    @(\"67-63=4\", \"93+4=97\", \"7+4=11\", \"53-51=2\", \"0+38=38\"),\n    @(\"9-1=8\", \"99-93=6\", \"63-34=29\", \"74-16=58\", \"18+27=45\"),\n    @(\"39+59=98\", \"37+25=62\", \"40+23=63\", \"53-2=51\", \"62-29=33\"),\n    @(\"73-8=65\", \"60-25=35\", \"0+1=1\", \"64+17=81\", \"10+59=69\"),\n    @(\"57-31=26\", \"38-0=38\", \"42+11=53\", \"52-28=24\", \"2+63=65\"),\n    @(\"77+21=98\", \"27-17=10\", \"76-45=31\", \"76-38=38\", \"77-41=36\")\n)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nif ($rowCount -ne $newValues.Count -or $colCount -ne $newValues[0].Count) {\n    throw \"Table shape ($rowCount x $colCount) does not match the replacement grid ($($newValues.Count) x $($newValues[0].Count)).\"\n}\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r-1][$c-1]\n    }\n}\n"}
